# maj template comment à la fin
#
# The "Comment" column (originally column J, the first of the
# Comment/SamplePortion/SamplePortionUnit/ResultUnit/NdfResult/AdfResult/
# AdlResult block) moves to the end of that block. In other words, every
# metadata row for columns J:P (field name, french label, #type tag,
# format description, example value) shifts one column to the left, with
# the old "Comment" column's content wrapping around into the last column
# (P) of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: field (API) names
$ws.Range("J1").Value = "SamplePortion"
$ws.Range("K1").Value = "SamplePortionUnit"
$ws.Range("L1").Value = "ResultUnit"
$ws.Range("M1").Value = "NdfResult"
$ws.Range("N1").Value = "AdfResult"
$ws.Range("O1").Value = "AdlResult"
$ws.Range("P1").Value = "Comment"

# Row 2: French labels
$ws.Range("J2").Value = "# Prise d'essai"
$ws.Range("K2").Value = "# Unité de mesure de la prise d’essai"
$ws.Range("L2").Value = "# Unité du résultat"
$ws.Range("M2").Value = "# Résultat NDF"
$ws.Range("N2").Value = "# Résultat ADF"
$ws.Range("O2").Value = "# Résultat ADL"
$ws.Range("P2").Value = "# Commentaire"

# Row 3: "#type" tags
$ws.Range("J3").Value = "#float"
$ws.Range("K3").Value = "#string"
$ws.Range("L3").Value = "#string"
$ws.Range("M3").Value = "#float"
$ws.Range("N3").Value = "#float"
$ws.Range("O3").Value = "#float"
$ws.Range("P3").Value = "#string"

# Row 4: format description
$ws.Range("J4").Value = "# format: nombre décimal, ne pas spécifier d'unité"
$ws.Range("K4").Value = "# format: texte"
$ws.Range("L4").Value = "# format: texte"
$ws.Range("M4").Value = "# format: nombre décimal ou NA"
$ws.Range("N4").Value = "# format: nombre décimal ou NA"
$ws.Range("O4").Value = "# format: nombre décimal ou NA"
$ws.Range("P4").Value = "# format: texte libre"

# Row 5: example value
$ws.Range("J5").Value = "# ex: 2.5"
$ws.Range("K5").Value = "# ex: mg"
$ws.Range("L5").Value = "# ex: mg/ml"
$ws.Range("M5").Value = "# 409.935 ou NA"
$ws.Range("N5").Value = "# 409.935 ou NA"
$ws.Range("O5").Value = "# 409.935 ou NA"
$ws.Range("P5").Value = ""
